$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.297536492347717
$ws.Range("B1").Value = 3.58649730682373
$ws.Range("C1").Value = 3.45299768447876
$ws.Range("D1").Value = 0.9610721468925476
$ws.Range("E1").Value = 1.142241597175598
